$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "43.813.30"
$ws.Cells.Item(2, 5).Value = "  -0.96%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.234.37"
$ws.Cells.Item(3, 5).Value = "  -2.08%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "314.16"
$ws.Cells.Item(5, 5).Value = "  -1.95%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "98.51"
$ws.Cells.Item(6, 5).Value = "  -5.45%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -3.21%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.25%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.531"
$ws.Cells.Item(9, 5).Value = "  -7.44%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "35.62"
$ws.Cells.Item(10, 5).Value = "  -8.80%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -2.65%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "7.33"
$ws.Cells.Item(12, 5).Value = "  -7.44%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -3.01%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "2.573.45"
$ws.Cells.Item(14, 5).Value = "  -2.03%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "2.232.10"
$ws.Cells.Item(15, 5).Value = "  -2.29%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.835"
$ws.Cells.Item(16, 5).Value = "  -5.42%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "13.90"
$ws.Cells.Item(17, 5).Value = "  -5.10%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "43.672.30"
$ws.Cells.Item(18, 5).Value = "  -1.05%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "13.02"
$ws.Cells.Item(19, 5).Value = "  -8.02%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0₃0962"
$ws.Cells.Item(20, 5).Value = "  -3.76%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.27"
$ws.Cells.Item(21, 5).Value = "  -4.94%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "65.77"
$ws.Cells.Item(22, 5).Value = "  -1.03%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "235.84"
$ws.Cells.Item(23, 5).Value = "  -1.01%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -7.73%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -8.52%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.40%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.06"
$ws.Cells.Item(27, 5).Value = "  -2.05%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -3.17%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "36.45"
$ws.Cells.Item(29, 5).Value = "  -7.20%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "5.95"
$ws.Cells.Item(30, 5).Value = "  -8.82%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "EthereumClassic"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "19.91"
$ws.Cells.Item(31, 5).Value = "  -3.17%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "Monero"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "156.75"
$ws.Cells.Item(32, 5).Value = "  -4.21%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0828"
$ws.Cells.Item(33, 5).Value = "  -6.37%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.28"
$ws.Cells.Item(34, 5).Value = "  -0.27%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.64"
$ws.Cells.Item(35, 5).Value = "  -3.09%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -9.02%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -3.24%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -3.55%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "15.53"
$ws.Cells.Item(39, 5).Value = "  -0.87%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -10.37%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.99"
$ws.Cells.Item(41, 5).Value = "  -11.92%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0305"
$ws.Cells.Item(42, 5).Value = "  -6.67%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +0.24%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.706.38"
$ws.Cells.Item(44, 5).Value = "  -3.81%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "81.99"
$ws.Cells.Item(45, 5).Value = "  -4.30%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -7.19%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "5.12"
$ws.Cells.Item(47, 5).Value = "  -5.33%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "101.26"
$ws.Cells.Item(48, 5).Value = "  -3.46%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "71.02"
$ws.Cells.Item(49, 5).Value = "  -6.00%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "55.96"
$ws.Cells.Item(50, 5).Value = "  -6.55%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -0.75%  "
